# Facutlty Testng file modification
# Update the sample Faculty row's address fields (State/Street/City/PostalCode).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Faculty")

$ws.Range("B2").Value = "Montana"
$ws.Range("C2").Value = "7244 Rau Station"
$ws.Range("D2").Value = "Russton"
$ws.Range("E2").Value = "65762-0220"
